$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row (at row 2), shifting
# existing data down.
$ws.Rows.Item(2).Resize(2).Insert() | Out-Null

# The insert copies formatting from the row above (the bold/bordered
# header), so reset the new rows back to the plain default style used
# by the other data rows.
$ws.Rows.Item(2).Resize(2).ClearFormats() | Out-Null

# Populate the two newly inserted rows with the new data values.
$ws.Range("A2").Value = 0.0937678143382072
$ws.Range("B2").Value = 0.0774271711707115
$ws.Range("C2").Value = -0.1965458989143371

$ws.Range("A3").Value = -0.1837177276611328
$ws.Range("B3").Value = 0.2562579810619354
$ws.Range("C3").Value = 0.0125227374956011

# Remove the trailing three old rows (now shifted to rows 22-24) so the
# data range ends at row 21 (A1:C21), matching the target dimension.
$ws.Rows.Item(22).Resize(3).Delete() | Out-Null
